# Petty Cash Book update - 14-Jun-2021, midday update.
# This script clears out the transaction entries recorded for 7-12 Jun 2021
# (rows 3 through 42 of "Sheet1") and starts a fresh entry for 14-Jun-2021,
# carrying forward the new opening balance.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New opening balance (SALDO AWAL) for this update.
$ws.Range("E2").Value = 276925

# Row 3 keeps its date (now 14-Jun-2021 = 44361) and its "Wages Expense"
# label, but the Debit entry that used to sit in D3 is removed.
$ws.Range("A3").Value = 44361
$ws.Range("D3").Clear()

# Rows 4 through 42 had all of their Tgl/Keterangan/Debit/Credit entries
# for the now-superseded days (7-12 Jun 2021); clear them out entirely so
# only the running "Saldo" formula in column E remains.
$ws.Range("A4:D42").Clear()

# Restore the usual cursor position (top of the sheet, cell D4) instead of
# where it had scrolled to previously.
$ws.Activate()
$ws.Range("D4").Select()
